$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Replace("15", "16")
$ws.Range("C9").Replace("4/8/2024", "4/15/2024")
$ws.Range("C9").Replace("4/14/2024", "4/21/2024")

# --- Cell type changes: number -> placeholder text "0" (shared style 14, string "0") ---
# Source: D14 (style 14, text "0"), untouched by this edit

# --- Type-change cells (style must flip between text-placeholder and numeric) ---
$ws.Range("D14").Copy($ws.Range("C17"))
$ws.Range("F14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 3
$ws.Range("H14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -33.333333333333
$ws.Range("D14").Copy($ws.Range("C20"))
$ws.Range("F14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("F14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = 1
$ws.Range("D14").Copy($ws.Range("D33"))
$ws.Range("E14").Copy($ws.Range("E33"))

# --- Simple numeric value updates ---
$ws.Range("N15").Value = -90
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 73.333333333333
$ws.Range("M16").Value = -3.703703703703
$ws.Range("N16").Value = -86.170212765957
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 2
$ws.Range("H17").Value = -84.615384615384
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = -50
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -68.354430379746
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = -23.255813953488
$ws.Range("L18").Value = 50
$ws.Range("M18").Value = -51.470588235294
$ws.Range("N18").Value = -90.934065934065
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 15.625
$ws.Range("I19").Value = 131
$ws.Range("J19").Value = 127
$ws.Range("K19").Value = 3.149606299212
$ws.Range("L19").Value = -23.391812865497
$ws.Range("M19").Value = 42.391304347826
$ws.Range("N19").Value = -11.486486486486
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -18.181818181818
$ws.Range("J20").Value = 33
$ws.Range("K20").Value = 72.727272727272
$ws.Range("L20").Value = 83.870967741935
$ws.Range("M20").Value = 5.555555555555
$ws.Range("N20").Value = -90.671031096563
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -27.777777777777
$ws.Range("F21").Value = 65
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = -10.958904109589
$ws.Range("I21").Value = 274
$ws.Range("J21").Value = 286
$ws.Range("K21").Value = -4.195804195804
$ws.Range("L21").Value = -1.083032490974
$ws.Range("M21").Value = 1.481481481481
$ws.Range("N21").Value = -80.442541042112
$ws.Range("M22").Value = 150
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -23.684210526315
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 112
$ws.Range("H24").Value = 1.785714285714
$ws.Range("I24").Value = 412
$ws.Range("J24").Value = 489
$ws.Range("K24").Value = -15.746421267893
$ws.Range("L24").Value = -28.842832469775
$ws.Range("M24").Value = 19.07514450867
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 5.882352941176
$ws.Range("F25").Value = 73
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = 52.083333333333
$ws.Range("I25").Value = 256
$ws.Range("J25").Value = 248
$ws.Range("K25").Value = 3.225806451612
$ws.Range("L25").Value = -28.690807799442
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -19.230769230769
$ws.Range("I26").Value = 95
$ws.Range("J26").Value = 114
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = -8.653846153846
$ws.Range("M26").Value = -12.037037037037
$ws.Range("I28").Value = 18
$ws.Range("K28").Value = 260
$ws.Range("L28").Value = 80
$ws.Range("F31").Value = 1
